$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild rows 3-13 of the ITT table in the new order described by the commit:
# the "Male condom attitudes index" / "Contraceptive knowledge index" blocks move
# up, directly under their related outcome, and every row picks up its own
# 95% CI columns (crude + adjusted).

# Row 3: Male condom attitudes index
$ws.Range("A3").Value = "'Male condom attitudes index"
$ws.Range("B3").Value = "'0.07 (1.03)"
$ws.Range("C3").Value = "'0.04 (0.97)"
$ws.Range("D3").Value = "'   -0.031"
$ws.Range("E3").Value = "'(-0.167, 0.106)"
$ws.Range("F3").Value = "'"
$ws.Range("G3").Value = "'"
$ws.Range("H3").Value = "'    0.011"
$ws.Range("I3").Value = "'(-0.103, 0.125)"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "'"

# Row 4: Used male condom at most recent sex (0/1)
$ws.Range("A4").Value = "'Used male condom at most recent sex (0/1)"
$ws.Range("B4").Value = "'528 (43.07%)"
$ws.Range("C4").Value = "'540 (44.85%)"
$ws.Range("D4").Value = "'    0.018"
$ws.Range("E4").Value = "'(-0.065, 0.101)"
$ws.Range("F4").Value = "'     1.08"
$ws.Range("G4").Value = "'(0.77, 1.50)"
$ws.Range("H4").Value = "'   -0.002"
$ws.Range("I4").Value = "'(-0.046, 0.041)"
$ws.Range("J4").Value = "'     0.98"
$ws.Range("K4").Value = "'(0.77, 1.23)"

# Row 5: Contraceptive knowledge index
$ws.Range("A5").Value = "'Contraceptive knowledge index"
$ws.Range("B5").Value = "'0.34 (0.98)"
$ws.Range("C5").Value = "'0.37 (0.94)"
$ws.Range("D5").Value = "'    0.029"
$ws.Range("E5").Value = "'(-0.090, 0.148)"
$ws.Range("F5").Value = "'"
$ws.Range("G5").Value = "'"
$ws.Range("H5").Value = "'    0.029"
$ws.Range("I5").Value = "'(-0.068, 0.125)"
$ws.Range("J5").Value = "'"
$ws.Range("K5").Value = "'"

# Row 6: Modern contraceptive methods known (n)
$ws.Range("A6").Value = "'Modern contraceptive methods known (n)"
$ws.Range("B6").Value = "'5.87 (2.52)"
$ws.Range("C6").Value = "'5.84 (2.58)"
$ws.Range("D6").Value = "'   -0.033"
$ws.Range("E6").Value = "'(-0.332, 0.265)"
$ws.Range("F6").Value = "'"
$ws.Range("G6").Value = "'"
$ws.Range("H6").Value = "'   -0.020"
$ws.Range("I6").Value = "'(-0.327, 0.287)"
$ws.Range("J6").Value = "'"
$ws.Range("K6").Value = "'"

# Row 7: Discussed contraceptive use with recent partner (0/1)
$ws.Range("A7").Value = "'Discussed contraceptive use with recent partner (0/1)"
$ws.Range("B7").Value = "'890 (72.89%)"
$ws.Range("C7").Value = "'881 (73.72%)"
$ws.Range("D7").Value = "'    0.008"
$ws.Range("E7").Value = "'(-0.056, 0.072)"
$ws.Range("F7").Value = "'     1.04"
$ws.Range("G7").Value = "'(0.75, 1.45)"
$ws.Range("H7").Value = "'   -0.012"
$ws.Range("I7").Value = "'(-0.065, 0.040)"
$ws.Range("J7").Value = "'     0.94"
$ws.Range("K7").Value = "'(0.72, 1.24)"

# Row 8: Can identify a female condom (0/1)
$ws.Range("A8").Value = "'Can identify a female condom (0/1)"
$ws.Range("B8").Value = "'678 (55.30%)"
$ws.Range("C8").Value = "'724 (60.13%)"
$ws.Range("D8").Value = "'    0.048"
$ws.Range("E8").Value = "'(-0.014, 0.111)"
$ws.Range("F8").Value = "'     1.22"
$ws.Range("G8").Value = "'(0.94, 1.57)"
$ws.Range("H8").Value = "'    0.048+"
$ws.Range("I8").Value = "'(0.001, 0.096)"
$ws.Range("J8").Value = "'     1.23*"
$ws.Range("K8").Value = "'(1.00, 1.52)"

# Row 9: Would be willing to try a female condom (0/1)
$ws.Range("A9").Value = "'Would be willing to try a female condom (0/1)"
$ws.Range("B9").Value = "'714 (58.24%)"
$ws.Range("C9").Value = "'715 (59.39%)"
$ws.Range("D9").Value = "'    0.011"
$ws.Range("E9").Value = "'(-0.037, 0.060)"
$ws.Range("F9").Value = "'     1.05"
$ws.Range("G9").Value = "'(0.86, 1.28)"
$ws.Range("H9").Value = "'   -0.006"
$ws.Range("I9").Value = "'(-0.054, 0.043)"
$ws.Range("J9").Value = "'     0.98"
$ws.Range("K9").Value = "'(0.80, 1.20)"

# Row 10: Female condom attitudes index
$ws.Range("A10").Value = "'Female condom attitudes index"
$ws.Range("B10").Value = "'0.01 (0.88)"
$ws.Range("C10").Value = "'-0.02 (0.86)"
$ws.Range("D10").Value = "'   -0.031"
$ws.Range("E10").Value = "'(-0.155, 0.093)"
$ws.Range("F10").Value = "'"
$ws.Range("G10").Value = "'"
$ws.Range("H10").Value = "'   -0.000"
$ws.Range("I10").Value = "'(-0.111, 0.111)"
$ws.Range("J10").Value = "'"
$ws.Range("K10").Value = "'"

# Row 11: Has ever used a female condom (0/1)
$ws.Range("A11").Value = "'Has ever used a female condom (0/1)"
$ws.Range("B11").Value = "'71 (5.79%)"
$ws.Range("C11").Value = "'86 (7.14%)"
$ws.Range("D11").Value = "'    0.014"
$ws.Range("E11").Value = "'(-0.005, 0.032)"
$ws.Range("F11").Value = "'     1.25"
$ws.Range("G11").Value = "'(0.92, 1.71)"
$ws.Range("H11").Value = "'    0.013"
$ws.Range("I11").Value = "'(-0.005, 0.032)"
$ws.Range("J11").Value = "'     1.25"
$ws.Range("K11").Value = "'(0.91, 1.73)"

# Row 12: Used a female condom in last 6 months (0/1)
$ws.Range("A12").Value = "'Used a female condom in last 6 months (0/1)"
$ws.Range("B12").Value = "'28 (2.28%)"
$ws.Range("C12").Value = "'34 (2.82%)"
$ws.Range("D12").Value = "'    0.005"
$ws.Range("E12").Value = "'(-0.006, 0.017)"
$ws.Range("F12").Value = "'     1.24"
$ws.Range("G12").Value = "'(0.78, 1.98)"
$ws.Range("H12").Value = "'    0.007"
$ws.Range("I12").Value = "'(-0.004, 0.018)"
$ws.Range("J12").Value = "'     1.33"
$ws.Range("K12").Value = "'(0.86, 2.05)"

# Row 13: Used a female condom at most recent sex (0/1)
$ws.Range("A13").Value = "'Used a female condom at most recent sex (0/1)"
$ws.Range("B13").Value = "'8 (0.65%)"
$ws.Range("C13").Value = "'9 (0.75%)"
$ws.Range("D13").Value = "'    0.001"
$ws.Range("E13").Value = "'(-0.006, 0.008)"
$ws.Range("F13").Value = "'     1.15"
$ws.Range("G13").Value = "'(0.43, 3.09)"
$ws.Range("H13").Value = "'    0.005"
$ws.Range("I13").Value = "'(-0.003, 0.012)"
$ws.Range("J13").Value = "'     2.02"
$ws.Range("K13").Value = "'(0.66, 6.22)"
